$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F9").Value = 1404
    $ws.Range("F22").Value = 2567
    $ws.Range("F23").Value = 37
}
